$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing rows 15..39 down to 16..40,
# then fill the new row 15 with the weekly data point.
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = Get-Date -Year 2021 -Month 9 -Day 28 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 100112012
$ws.Range("G15").Value = "Espinaca"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 250
$ws.Range("K15").Value = 800
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = 850
$ws.Range("N15").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 283
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = "Hortaliza"
